# Populate the workbook's "Sheet1" into a proper "Customers" sheet by
# inserting a new leading "Name" column (done while working with
# openpyxl-populated data in the real edit; reproduced here via COM).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rename Sheet1 -> Customers
$ws.Name = "Customers"

# Insert a new column A (shifts First Name, Last Name, ... one column right)
$ws.Columns.Item(1).Insert()

# New header for the inserted column
$ws.Range("A1").Value = "Name"

# Match the column's best-fit width for the short "Name" header
$ws.Columns.Item(1).ColumnWidth = 5.3
